$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 0.0207
$ws.Range("E2").Value = -0.0005999999999999964
$ws.Range("F2").Value = 0.0665
$ws.Range("G2").Value = 0.09108869671915629
$ws.Range("H2").Value = 0.09108869671915629
$ws.Range("I2").Value = 0.07494531258736331
$ws.Range("J2").Value = 0.0577455438301705
$ws.Range("K2").Value = 1070.6
$ws.Range("L2").Value = 0.03909696456221333
$ws.Range("M2").Value = 614.6110000000001
$ws.Range("N2").Value = 0.06056474182104849
$ws.Range("O2").Value = 0.5740808892209978
$ws.Range("P2").Value = 610.0810000000001
$ws.Range("Q2").Value = 0.06011834844304298
$ws.Range("R2").Value = 0.5698496170371756
$ws.Range("S2").Value = 4.530000000000001
$ws.Range("T2").Value = 0.007370515659498447
$ws.Range("U2").Value = 4760.9
$ws.Range("V2").Value = 0.4691466298778084
$ws.Range("W2").Value = 0.08552891756269879
$ws.Range("X2").Value = 0.06176853825336105
$ws.Range("Y2").Value = 0.02376037930933774
$ws.Range("Z2").Value = 2.714301799749581
$ws.Range("AA2").Value = 0.1909303960774392
$ws.Range("AB2").Value = 0.0519107813444111
$ws.Range("AC2").Value = 0.139019614733028
$ws.Range("AD2").Value = 3840.2
$ws.Range("AE2").Value = 294.2875817885643
$ws.Range("AF2").Value = 4134.487581788564
$ws.Range("AG2").Value = -626.4124182114356
$ws.Range("AH2").Value = 0.2894795152533794
$ws.Range("AI2").Value = 0.211125419549952
$ws.Range("AJ2").Value = -0.06578865266224525
$ws.Range("AK2").Value = -0.04226180614399212
$ws.Range("AL2").Value = 100.6
$ws.Range("AM2").Value = 100.6
$ws.Range("AN2").Value = 1.6294819026605
$ws.Range("AO2").Value = 20.3648111332008
$ws.Range("AP2").Value = -0.2658006611836193
$ws.Range("AQ2").Value = 20.3648111332008

# Row 3 updates
$ws.Range("D3").Value = 0.0621
$ws.Range("E3").Value = 0.0469
$ws.Range("G3").Value = 0.1481522144880437
$ws.Range("H3").Value = 0.1481522144880437
$ws.Range("I3").Value = 0.1162249401638084
$ws.Range("J3").Value = 0.08841801773210803
$ws.Range("K3").Value = 369.9
$ws.Range("L3").Value = 0.08122351287850507
$ws.Range("M3").Value = 96.23
$ws.Range("N3").Value = 0.02287867620836404
$ws.Range("O3").Value = 0.2601513922681806
$ws.Range("P3").Value = 91.7
$ws.Range("Q3").Value = 0.02180166900454102
$ws.Range("R3").Value = 0.2479048391457151
$ws.Range("S3").Value = 4.530000000000001
$ws.Range("T3").Value = 0.04707471682427519
$ws.Range("U3").Value = 1635
$ws.Range("V3").Value = 0.3887211431016856
$ws.Range("W3").Value = 0.09952644890491309
$ws.Range("X3").Value = 0.05238534091293877
$ws.Range("Y3").Value = 0.04714110799197432
$ws.Range("Z3").Value = 2.721790580922783
$ws.Range("AA3").Value = 0.2406553278471152
$ws.Range("AB3").Value = 0.05098925142216552
$ws.Range("AC3").Value = 0.1896660764249497
$ws.Range("AD3").Value = 225.8
$ws.Range("AF3").Value = 225.8
$ws.Range("AG3").Value = -1409.2
$ws.Range("AH3").Value = 0.05094880299645749
$ws.Range("AI3").Value = 0.05114962056858081
$ws.Range("AJ3").Value = -0.5038435410633201
$ws.Range("AK3").Value = -0.5069976614499011
$ws.Range("AN3").Value = 0.3930374238468233
$ws.Range("AP3").Value = -2.452915578764143

# Row 3 cell deletion (F3 removed, shifting cells left not required - clear content only)
$ws.Range("F3").ClearContents()

# Row 4 updates
$ws.Range("D4").Value = -0.0207
$ws.Range("E4").Value = -0.0481
$ws.Range("F4").Value = 0.0665
$ws.Range("G4").Value = 0.07970528842573733
$ws.Range("H4").Value = 0.07970528842573733
$ws.Range("I4").Value = 0.06671057920120756
$ws.Range("J4").Value = 0.05205131127682683
$ws.Range("K4").Value = 700.7
$ws.Range("L4").Value = 0.03069328182013308
$ws.Range("M4").Value = 518.3810000000001
$ws.Range("N4").Value = 0.08724162304986623
$ws.Range("O4").Value = 0.7398044812330528
$ws.Range("P4").Value = 518.3810000000001
$ws.Range("Q4").Value = 0.08724162304986623
$ws.Range("R4").Value = 0.7398044812330528
$ws.Range("U4").Value = 3125.9
$ws.Range("V4").Value = 0.5260775172924486
$ws.Range("W4").Value = 0.07153138622048449
$ws.Range("X4").Value = 0.07115173559378332
$ws.Range("Y4").Value = 0.0003796506267011723
$ws.Range("Z4").Value = 2.712812815738372
$ws.Range("AA4").Value = 0.1412054643077631
$ws.Range("AB4").Value = 0.05283231126665668
$ws.Range("AC4").Value = 0.08837315304110638
$ws.Range("AD4").Value = 3614.4
$ws.Range("AE4").Value = 294.2875817885643
$ws.Range("AF4").Value = 3908.687581788564
$ws.Range("AG4").Value = 782.7875817885642
$ws.Range("AH4").Value = 0.3967974041482373
$ws.Range("AI4").Value = 0.257683028212946
$ws.Range("AJ4").Value = 0.1164050481554663
$ws.Range("AK4").Value = 0.06500107027374248
$ws.Range("AL4").Value = 100.6
$ws.Range("AM4").Value = 100.6
$ws.Range("AN4").Value = 2.028055212658512
$ws.Range("AO4").Value = 15.10337972166998
$ws.Range("AP4").Value = 0.4392254414704097
$ws.Range("AQ4").Value = 15.10337972166998
